$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (I1) and IF (J1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style/formatting used by the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# I0 / IF values per data row (row -> I, J)
$data = @{
    2  = @(1, 7)
    3  = @(1, 6)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(1, 7)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 9)
    10 = @(1, 8)
    11 = @(1, 6)
    12 = @(1, 8)
    13 = @(1, 7)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 8)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 7)
    21 = @(3, 7)
    22 = @(1, 6)
    23 = @(1, 6)
    24 = @(1, 7)
    25 = @(1, 7)
    26 = @(1, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
